$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I12").Value = 0.2473204566274169
$ws.Range("J12").Value = 0.008750979825177335
$ws.Range("K12").Value = 0.8089196720461374
$ws.Range("L12").Value = 2.101860250809634

$ws.Range("I13").Value = 0.2666942839323744
$ws.Range("J13").Value = 0.008502164597911139
$ws.Range("K13").Value = 0.2680339375067249
$ws.Range("L13").Value = 1.532865872683723

$ws.Range("I14").Value = 0.3380914860526574
$ws.Range("J14").Value = 0.03257293951794873
$ws.Range("K14").Value = 1.009832860521988
$ws.Range("L14").Value = 2.815720396193842

$ws.Range("I15").Value = 0.3951692466145622
$ws.Range("J15").Value = 0.04746688694852878
$ws.Range("K15").Value = 0.569185250568909
$ws.Range("L15").Value = 1.739304231923645

$ws.Range("I17").Value = 0.4561555070866828
$ws.Range("J17").Value = 0.06756879384457612
$ws.Range("K17").Value = 0.4030054558325659
$ws.Range("L17").Value = 1.590884813525737
